$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data point (2026/02/26, 木, 14, 201) was recorded. It belongs
# right before the existing 2026/12/29 block, so insert a fresh row at 888;
# everything previously at/after row 888 shifts down by one
# (old row 888 -> 889, ..., old row 929 -> 930).
$ws.Rows.Item(888).Insert()

# Column A holds a date written as plain text (e.g. "2026/12/29"), not a
# real date value. Force text formatting before assigning so it isn't
# auto-converted to a date serial, then drop the format change so the
# cell ends up with no explicit style, matching its neighbours.
$ws.Cells.Item(888, 1).NumberFormat = "@"
$ws.Cells.Item(888, 1).Value = "2026/02/26"
$ws.Cells.Item(888, 1).ClearFormats()

$ws.Cells.Item(888, 2).Value = "木"
$ws.Cells.Item(888, 3).Value = 14
$ws.Cells.Item(888, 4).Value = 201
